# CDS Input file updates
# Replace the (old, simpler) "ParticipantsTab" Neo4j query in B2 with the new
# query that adds diagnosis/genomic_info lookups and sorts the collected
# sample ids via apoc.coll.sort(), then refresh the view/selection state to
# match how the workbook was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newParticipantQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['BAM']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

# B2 holds the "query" column for the ParticipantsTab row; only its text
# content changes (style/column position stay the same).
$ws.Range("B2").Value = $newParticipantQuery

# The new query text wraps across more lines, so the row grows taller.
$ws.Rows.Item(2).RowHeight = 279

# Window/selection state left behind after the edit.
$ws.Range("C4").Select()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 3
    $win.ScrollColumn = 1
}
